$d = $word.ActiveDocument

# 1. Collapse the split " " + "are" + " " runs (with grammar-check proofErr
#    markers around "are") into a single " are " run.
$d.Content.Find.Execute(" are ", $false, $false, $false, $false, $false,
                         $true, 1, $false, " are ", 2) | Out-Null

# 2. Append two empty paragraphs at the very end of the document body
#    (before the final section break).
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
